$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Growth_protocol")
$ws.Activate()
$ws.Range("G2:G11").Select() | Out-Null
$ws.Range("G2:G11").ClearContents()
